# Added Python and MySQL cheatsheet
# ------------------------------------------------------------------
# Helper: write a run-formatted ("rich text") value into a cell.
# $runs is an array of 2-element arrays: @(text, isBold)
# The very first run is only touched if it needs to be bold (matches
# Excel's own habit of leaving an un-styled leading run with no <rPr/>).
# ------------------------------------------------------------------
function Set-RichCell {
    param($cell, $runs)
    $full = ""
    foreach ($r in $runs) { $full += $r[0] }
    $cell.Value = $full
    $pos = 1
    $i = 0
    foreach ($r in $runs) {
        $txt = $r[0]
        $bold = $r[1]
        $len = $txt.Length
        if ($len -gt 0) {
            if ($i -eq 0) {
                if ($bold) {
                    $cell.Characters($pos, $len).Font.Bold = $true
                }
            } else {
                $cell.Characters($pos, $len).Font.Bold = $bold
            }
        }
        $pos += $len
        $i += 1
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# New cheatsheet rows: ENUM / AUTO_INCREMENT / DEFAULT (MySQL)
#
# NOTE: the order these distinct strings are first assigned controls
# the shared-string table index order on save, so the cells are
# populated in the same order the original author entered them
# (column A for all three rows, then column B, then column C) rather
# than strictly row-by-row.
# ------------------------------------------------------------------

# -- Column A (Concept) for the three new rows --------------------
Set-RichCell $ws.Range("A27") @(
    ,@("Choose from only a select number of items when creating the table with ", $false)
    ,@("ENUM", $true)
)

Set-RichCell $ws.Range("A28") @(
    ,@("Automatically auto-increment the Primary Key and add NOT NULL", $false)
)

Set-RichCell $ws.Range("A29") @(
    ,@("Set a default value for columns in a table", $false)
)

# -- Column B (Code) ------------------------------------------------
Set-RichCell $ws.Range("B28") @(
    ,@("CREATE TABLE Person(`r`n    personID INT PRIMARY KEY NOT NULL ", $false)
    ,@("AUTO_INCREMENT", $true)
    ,@(",`r`n    name VARCHAR(50) NOT NULL,`r`n    gender ENUM('M', 'F', 'O') NOT NULL,`r`n);", $false)
)

Set-RichCell $ws.Range("B27") @(
    ,@("CREATE TABLE Person(`r`n    personID INT PRIMARY KEY,`r`n    name VARCHAR(50),`r`n   ", $false)
    ,@(" gender ENUM('M', 'F', 'O'),", $true)
    ,@("`r`n);", $false)
)

Set-RichCell $ws.Range("B29") @(
    ,@("CREATE TABLE Score(`r`n    scoreID INT PRIMARY KEY NOT NULL AUTO_INCREMENT,`r`n    score1 INT ", $false)
    ,@("DEFAULT 0", $true)
    ,@(",`r`n    score2 INT ", $false)
    ,@("DEFAULT 0", $true)
    ,@(",`r`n)", $false)
)

# -- Column C (Explained) -------------------------------------------
Set-RichCell $ws.Range("C29") @(
    ,@("We can create default values with the ", $false)
    ,@("DEFAULT", $true)
    ,@(" keyword, then insert whatever we'd like afterwards.", $false)
)

Set-RichCell $ws.Range("C28") @(
    ,@("The ", $false)
    ,@("AUTO_INCREMENT", $true)
    ,@(" keyword makes it so that we don't need to manually insert the primary key all the time", $false)
)

# ------------------------------------------------------------------
# Row heights - Excel re-wrapped these rows once column A widened
# and the new rows were appended; reproduce the resulting heights.
# ------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 105
$ws.Rows.Item(13).RowHeight = 150
$ws.Rows.Item(17).RowHeight = 75
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 45
$ws.Rows.Item(24).RowHeight = 75
$ws.Rows.Item(26).RowHeight = 45
$ws.Rows.Item(27).RowHeight = 75
$ws.Rows.Item(28).RowHeight = 75
$ws.Rows.Item(29).RowHeight = 75

# ------------------------------------------------------------------
# Column A - widen, no longer a "best fit" width
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19

# ------------------------------------------------------------------
# Table1 now covers the three new rows
# ------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E29"))

# ------------------------------------------------------------------
# Selection / scroll position, as left by the author after editing
# ------------------------------------------------------------------
$ws.Range("D28").Select()
